$d = $word.ActiveDocument

# --- 1. Reformat paragraph 7 (1-indexed): "Модель стоимости жилья  Магнитогорске" ---
$p6 = $d.Paragraphs.Item(7)
$r6 = $p6.Range
$r6.MoveEnd(1, -1)
$r6.Delete()
$xml6 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Модель стоимости </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>жилья</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve"> Магнитогорске</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r6.InsertXML($xml6)

# --- 2. Reformat paragraph 16 (1-indexed): "Для трех квартир ... модели ни для контроля качества." ---
$p15 = $d.Paragraphs.Item(16)
$r15 = $p15.Range
$r15.MoveEnd(1, -1)
$r15.Delete()
$xml15 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Для трех квартир цена не указана, поэтому удалим эти строки, так как они не подходят ни для </w:t></w:r><w:r><w:t xml:space="preserve">обучения </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>модели</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> ни для контроля качества.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r15.InsertXML($xml15)

# --- 3. Append new paragraphs after the last paragraph ---
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$rLast = $last.Range
$rLast.InsertParagraphAfter()
$newp = $d.Paragraphs.Item($d.Paragraphs.Count)
$rNew = $newp.Range
$rNew.Collapse(0)
$xmlNew = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Посмотрим количество пропусков в столбце с данными о районе.</w:t></w:r></w:p><w:p><w:r><w:t>Заполним пропуски значением неизвестно.</w:t></w:r></w:p><w:p><w:r><w:t>Построение модели</w:t></w:r></w:p><w:p><w:r><w:t>Подготовка обучающей и валидационной выборки.</w:t></w:r></w:p><w:p><w:r><w:t>Выделим характеристики, которые будем использовать для обучения модели.</w:t></w:r></w:p><w:p><w:r><w:t>не влияет на цену квартиры</w:t></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:t>может оказывать влияние</w:t></w:r></w:p><w:p><w:r><w:t>вместо него будем использовать более информативные столбцы, сгенерированные на основе данных из этого столбца</w:t></w:r></w:p><w:p><w:r><w:t>целевой признак</w:t></w:r></w:p><w:p><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>id</w:t></w:r><w:r><w:t xml:space="preserve"> для демонстрации на собеседовании</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Разобьем выборки на обучающую и валидационную в отношении </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>4 :</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> 1.</w:t></w:r></w:p><w:p><w:r><w:t>Предварительна обработка данных</w:t></w:r></w:p><w:p><w:r><w:t>Выделим категориальные и числовые признаки.</w:t></w:r></w:p><w:p><w:r><w:t>Количество комнат попадает в категориальные признаки, потому что содержит значение «многоквартирная», которое нельзя заменить конкретным числовым значением.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Преобразуем категориальные признаки в числа с помощью порядкового кодирования. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rNew.InsertXML($xmlNew)

# clean up the trailing empty paragraph mark left over from InsertParagraphAfter
$trailing = $d.Paragraphs.Item($d.Paragraphs.Count)
$trStart = $trailing.Range.Start
$docEnd = $d.Content.End
$delRange = $d.Range($trStart - 1, $docEnd)
$delRange.Delete()

Write-Host "Final paragraph count:" $d.Paragraphs.Count
